$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (week of 2023-12-15) arrived and needs to be
# inserted ahead of the existing row 46, pushing rows 46:52 down to 47:53.
$ws.Rows("46:46").Insert()

$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C46").Value = "Ñuble"
$ws.Range("D46").Value = 45275
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = 100112039
$ws.Range("G46").Value = "Ciboulette"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 150
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = 2000
$ws.Range("N46").Value = "$/docena de atados"
$ws.Range("O46").Value = "Región Metropolitana"
$ws.Range("P46").Value = 667
$ws.Range("Q46").Value = 3
$ws.Range("R46").Value = "Hortaliza"
